$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# New header cells: Wins / Losses / Ties in AD1:AF1, matching the
# bold/bordered/centered style already used by the other header cells.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial($xlPasteFormats)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record (Wins=49, Losses=63, Ties=0) for every
# data row (2 through 46).
for ($row = 2; $row -le 46; $row++) {
    $ws.Cells.Item($row, 30).Value = 49
    $ws.Cells.Item($row, 31).Value = 63
    $ws.Cells.Item($row, 32).Value = 0
}

"done"
